# regen save_data to use K (strikeouts) instead of Strike# (TB-derived) for column G,
# then recompute the dependent std/mean/s_vals style summary stats.
# The new K values come from the refreshed source box-score data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value (column G), pulled from the regenerated save_data source.
$kValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 1
    6 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 2
    14 = 4
    15 = 0
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 0
    24 = 3
    25 = 0
    26 = 1
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 0
    32 = 0
    33 = 3
    34 = 3
    35 = 0
    36 = 1
    37 = 0
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 3
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    49 = 1
    50 = 0
    51 = 2
    52 = 3
    53 = 1
    54 = 2
    55 = 1
    56 = 5
    57 = 0
    58 = 1
    59 = 0
    60 = 2
    61 = 2
    62 = 1
    63 = 1
    64 = 1
    65 = 1
    66 = 1
    67 = 0
    68 = 1
    69 = 2
    70 = 2
    71 = 3
    72 = 2
    74 = 2
    75 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
